# Generate Report for Handoff
# Adds a new row (for file bea98671-cec3-41e3-b5be-5a1c5c400b33) to the
# Overview, zh-cn and de-de sheets, mirroring the existing row that was
# generated for 9458f540-1bfc-4017-9c00-dea3fdc88b08.

$wb = $excel.ActiveWorkbook

$newGuid = "bea98671-cec3-41e3-b5be-5a1c5c400b33"
$newHash = "96d6ef1229bda4165f79294360227d3429a2f1fa"
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3fca21c6d5122741fc290e773d5ccc50d4aac52/e2e/"
$hyperlinkColor = 15570276   # BGR long for FF6495ED ("cornflower blue"), matches existing hyperlink style

# -----------------------------------------------------------------
# Sheet "Overview" -> new row 3
# -----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "$newGuid.md"

$wsOverview.Range("B3").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "$baseUrl$newGuid.md", "", "", "e2e\$newGuid.md")
$wsOverview.Range("B3").Font.Color = $hyperlinkColor
$wsOverview.Range("B3").Font.Underline = 2

$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = "'"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"

$wsOverview.Range("G3").Value = "2016-08-26 20:38:56"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# -----------------------------------------------------------------
# Sheet "zh-cn" -> new row 3
# -----------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = "$newGuid.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "$baseUrl$newGuid.md", "", "", "$newGuid.md")
$wsZhCn.Range("A3").Font.Color = $hyperlinkColor
$wsZhCn.Range("A3").Font.Underline = 2

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"

$wsZhCn.Range("G3").Value = "$newGuid.$newHash.zh-cn.xlf"

$wsZhCn.Range("H3").Value = "2016-08-26 20:38:51"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("I3").Value = "'"
$wsZhCn.Range("J3").Value = "'"

$wsZhCn.Range("K3").Value = "2016-08-26 20:38:35"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("L3").Value = "'"
$wsZhCn.Range("M3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = "'"

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# -----------------------------------------------------------------
# Sheet "de-de" -> new row 3
# -----------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = "$newGuid.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "$baseUrl$newGuid.md", "", "", "$newGuid.md")
$wsDeDe.Range("A3").Font.Color = $hyperlinkColor
$wsDeDe.Range("A3").Font.Underline = 2

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"

$wsDeDe.Range("G3").Value = "$newGuid.$newHash.de-de.xlf"

$wsDeDe.Range("H3").Value = "2016-08-26 20:38:56"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("I3").Value = "'"
$wsDeDe.Range("J3").Value = "'"

$wsDeDe.Range("K3").Value = "2016-08-26 20:38:35"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("L3").Value = "'"
$wsDeDe.Range("M3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("N3").Value = "'"
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = "'"

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
